$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells (row 1)
$ws.Range("D1").Value = "Y_min"
$ws.Range("E1").Value = "Y_max"
$ws.Range("F1").Value = "Y"

# Rename values in column K (A_eta_W -> W_aysm) and column L (eta -> Y) for rows 2-14
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 11).Value = "W_aysm"
    $ws.Cells.Item($r, 12).Value = "Y"
}

# Update the active selection to match the recorded view state
$ws.Range("K20").Select()
